$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 8 new daily rows (2021-09-02 .. 2021-09-09, Excel serials 44441-44448)
# right after the existing last data row (366), with zero counts, matching
# the existing table's data pattern.
$startRow = 367
$startSerial = 44441
$numNewRows = 8

for ($i = 0; $i -lt $numNewRows; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $startSerial + $i
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}

# Match the date-column formatting used by the rest of column A (bordered,
# bold, centered, custom date/time number format) by copying it from the
# last pre-existing row in that column.
$ws.Range("A366").Copy()
$ws.Range("A367:A374").PasteSpecial(-4122)
